# Update the cryptos list (Price and Volume(1h) columns) with freshly
# scraped values, as produced by the scheduled GitHub Actions scraper run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2; D="28.080.00"; E="  +1.28%  "}
    @{Row=3; D="1.790.17"; E="  +1.77%  "}
    @{Row=4; D="1.001"; E="  -0.50%  "}
    @{Row=5; D="323.41"; E="  -1.02%  "}
    @{Row=6; D="0.9997"; E="  -0.13%  "}
    @{Row=7; D="0.4313"; E="  -2.73%  "}
    @{Row=8; D="0.3624"; E="  -2.98%  "}
    @{Row=9; D="44.73"; E="  -2.15%  "}
    @{Row=10; D="0.07500"; E="  -3.57%  "}
    @{Row=11; D="1.113"; E="  -1.36%  "}
    @{Row=12; D="0.9998"; E="  -0.24%  "}
    @{Row=13; D="21.66"; E="  -0.80%  "}
    @{Row=14; D="6.146"; E=$null}
    @{Row=15; D=$null; E="  -0.86%  "}
    @{Row=16; D="1.788.85"; E="  +1.69%  "}
    @{Row=17; D="92.15"; E="  +0.43%  "}
    @{Row=18; D="0.00001064"; E="  -1.77%  "}
    @{Row=19; D="0.06350"; E="  +2.01%  "}
    @{Row=20; D="0.9999"; E="  -0.09%  "}
    @{Row=21; D="17.24"; E="  -1.27%  "}
    @{Row=22; D="5.964"; E="  -3.85%  "}
    @{Row=23; D="28.076.89"; E="  +1.14%  "}
    @{Row=24; D="11.38"; E="  -2.41%  "}
    @{Row=25; D="2.121"; E="  -8.72%  "}
    @{Row=26; D="158.85"; E="  +3.44%  "}
    @{Row=27; D="20.38"; E="  -2.37%  "}
    @{Row=28; D="1.993.60"; E="  +1.82%  "}
    @{Row=29; D=$null; E="  -8.71%  "}
    @{Row=30; D="127.18"; E="  -1.61%  "}
    @{Row=31; D="1.160"; E="  -4.70%  "}
    @{Row=32; D="5.721"; E="  -1.25%  "}
    @{Row=33; D="0.09006"; E="  -2.99%  "}
    @{Row=34; D="3.536"; E="  -3.68%  "}
    @{Row=35; D="12.59"; E="  -1.36%  "}
    @{Row=36; D="0.02319"; E="  -1.14%  "}
    @{Row=37; D="5.093"; E="  -0.37%  "}
    @{Row=38; D="0.6451"; E="  -1.40%  "}
    @{Row=39; D="0.2115"; E="  -3.88%  "}
    @{Row=40; D="0.06049"; E="  -1.65%  "}
    @{Row=41; D=$null; E="  -0.90%  "}
    @{Row=42; D="1.420"; E="  +0.29%  "}
    @{Row=43; D=$null; E="  -0.10%  "}
    @{Row=44; D="7.835"; E="  -2.70%  "}
    @{Row=45; D="13.64"; E="  -1.86%  "}
    @{Row=46; D="0.5982"; E="  -1.11%  "}
    @{Row=47; D="3.704"; E="  -1.55%  "}
    @{Row=48; D="124.72"; E="  -1.14%  "}
    @{Row=49; D=$null; E="  -1.17%  "}
    @{Row=50; D="1.152"; E="  +0.09%  "}
    @{Row=51; D="0.06945"; E="  +0.42%  "}
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($u.Row, 4)
        # Force text storage so number-like strings (e.g. "1.001",
        # "323.41") stay text the same way the scraped source data
        # always has (matches the original inlineStr cells), instead of
        # Excel auto-coercing them to numeric values. Reset back to the
        # default "Normal" style afterwards so no stray number format is
        # left behind on the cell.
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.Style = "Normal"
    }
    if ($null -ne $u.E) {
        $ws.Cells.Item($u.Row, 5).Value = $u.E
    }
}
